# "Megan push July 29" - update To Do List checkmarks/highlights and
# Analysis sheet findings text, then leave the Analysis tab active.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("To Do List")
$ws2 = $wb.Worksheets.Item("Analysis")

# ---------------------------------------------------------------
# Sheet 1: "To Do List"
# ---------------------------------------------------------------

# Megan's personal to-do column (E) - mark more items as done (strikethrough)
$ws1.Range("E5:E8").Font.Strikethrough = $true

# "Last Group Day" column (G) - mark two more items done (strikethrough),
# and add the new outstanding tasks at the bottom of the list.
$ws1.Range("G2").Font.Strikethrough = $true
$ws1.Range("G3").Font.Strikethrough = $true

$ws1.Range("G6").Value = "Bechdel Test - what is it"

$ws1.Range("G7").Value = "Bechdel pie charts"
$ws1.Range("G7").Font.Strikethrough = $true

$ws1.Range("G8").Value = "bechdel horizontal bar"
$ws1.Range("G8").Font.Strikethrough = $true

$ws1.Range("G9").Value = "Summarize Findings"
$ws1.Range("G10").Value = "General Conclusions in Slideshow"

# Highlight the newly-relevant outline row, and clear the highlight that
# is no longer needed.
$ws1.Range("A12").Interior.Color = 65535
$ws1.Range("A16").Interior.Pattern = -4142

# Update the selection on this sheet (it is no longer the active tab).
$ws1.Range("G8").Select()

# ---------------------------------------------------------------
# Sheet 2: "Analysis"
# ---------------------------------------------------------------

$ws2.Range("B7").Value = "Across the Internet Movie Data Base (IMDB), the ratings for top grossing films and best picture award winners are markedly similarly across the years with both pictures typically rating favorably.  Meanwhile, the Rotten Tomato ""Tomatometer"" shows a  greater variation between top grossing film scores and best picture film scores.  Particularly of interest are the films from 2014 (Transformers 2) vs. (Birdman).  Overall, best picture films score consistently favorably on the Tomatometer, whilst the more volatile score swings belong to the top grossing films."

$ws2.Rows.Item(7).RowHeight = 115.2

# Analysis becomes the active sheet / tab, with B7 selected.
$ws2.Activate()
$ws2.Range("B7").Select()
